# Training Dashboard refresh: restyle title/header fonts, swap data-row
# highlight colors based on recalculated validity, and bump the figures
# (period-to-expire / last-update / status) for the new "as of" date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Colors (VBA-style BGR long values) ---
$white = 16777215   # FFFFFF
$red   = 13551615   # FFC7CE (existing "NOT VALID" highlight fill)

# --- Title (A1): keep bold, drop the explicit 14pt size, make it white ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = $white

# --- Header row (A2:K2): keep bold, make the text white on the blue fill ---
$ws.Range("A2:K2").Font.Bold = $true
$ws.Range("A2:K2").Font.Color = $white

# --- Data rows: recompute PERIOD TO EXPIRE / LAST UPDATE / STATUS ---
$today = "16-Sep-2025"

$ws.Range("H3").Value = 20
$ws.Range("I3").Value = $today
$ws.Range("J3").Value = "NOT VALID"

$ws.Range("H4").Value = -42
$ws.Range("I4").Value = $today
$ws.Range("J4").Value = "NOT VALID"

$ws.Range("H5").Value = 204
$ws.Range("I5").Value = $today
$ws.Range("J5").Value = "VALID"

$ws.Range("H6").Value = 205
$ws.Range("I6").Value = $today
$ws.Range("J6").Value = "VALID"

$ws.Range("H7").Value = 304
$ws.Range("I7").Value = $today
$ws.Range("J7").Value = "VALID"

# --- Re-apply the NOT VALID / VALID row highlight now that the statuses
#     above have been recalculated (row 3 flips to NOT VALID => red) ---
$ws.Range("A3:K3").Interior.Color = $red
$ws.Range("A4:K4").Interior.Color = $red
$ws.Range("A5:K5").Interior.Color = $white
$ws.Range("A6:K6").Interior.Color = $white
$ws.Range("A7:K7").Interior.Color = $white
